$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "X" marker in row 3 from the "In Progress" column (B) to the
# "Done" column (D): the task is now complete.
$ws.Range("B3").Value = $null
$ws.Range("D3").Value = "X"

# Update the active cell selection shown in the sheet view.
$ws.Range("A5").Select()
